# Apply LOM3217.docx edit: insert manual line breaks (<w:br/>) within three
# paragraphs (Programa PT, Programa EN, Bibliografia), splitting each long run
# of text into several runs at specific sentence boundaries, per the commit diff.

$d = $word.ActiveDocument

function Replace-WithLineBreaks($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
    return $ok
}

# 1) Programa (Portuguese)
Replace-WithLineBreaks 'Sistemas macroscópicos e microscópicos. Postulados da termodinâmica. Equação fundamental. Equações de estado. Equação de Euler. Relação de Gibbs-Duhem. Equilíbrio termodinâmico. Derivadas Termodinâmicas. Potenciais Termodinâmicos: Helmholtz, Gibbs, Grande Canônico e Entalpia. Relações de Maxwell. Diagrama de Born. Redução de derivadas termodinâmicas. Formalismo microcanônico. Equação de Boltzmann. Exemplo: Modelo de Einstein de sólido cristalino. Formalismo canônico. Exemplo: Gás ideal clássico. Distribuição de Maxwell-Boltzmann. Formalismo grande canônico. Gases quânticos: férmions e bósons. Estatística de Bose-Einstein. Estatística de Fermi-Dirac. Exemplos: gás de elétrons e gás de fótons. Estatística de Planck.Aplicações: supercondutividade. Gás de elétrons em semicondutores. superfluidez no hélio líquido.' 'Sistemas macroscópicos e microscópicos. Postulados da termodinâmica. Equação fundamental. Equações de estado. Equação de Euler. Relação de Gibbs-Duhem. Equilíbrio termodinâmico. Derivadas Termodinâmicas. Potenciais Termodinâmicos: Helmholtz, Gibbs, Grande Canônico e Entalpia. Relações de Maxwell. Diagrama de Born. Redução de derivadas termodinâmicas. ^lFormalismo microcanônico. Equação de Boltzmann. Exemplo: Modelo de Einstein de sólido cristalino. Formalismo canônico. Exemplo: Gás ideal clássico. Distribuição de Maxwell-Boltzmann. Formalismo grande canônico. Gases quânticos: férmions e bósons. Estatística de Bose-Einstein. Estatística de Fermi-Dirac. Exemplos: gás de elétrons e gás de fótons. Estatística de Planck.^lAplicações: supercondutividade. Gás de elétrons em semicondutores. superfluidez no hélio líquido.' | Out-Null

# 2) Programa (English / italic)
Replace-WithLineBreaks 'Macroscopic and microscopic systems. Postulates of thermodynamics. Fundamental equation. Equations of state. Thermodynamic equilibrium. Thermodynamic derivatives. Thermodynamic potentials. Maxwell relations. Born diagram. Reduction of thermodynamic derivatives.Microcanonical formalism. Boltzmann equation. Einstein model of a crystalline solid. Canonical formalism. Example: ideal classical gas. Maxwell-Boltzmann distribution. Grand canonical formalism. Quantum gases. Fermions and bosons. Bose-Einstein distribution. Fermi-Dirac distribution. Examples: electron gas and photon gas. Planck distribution.Applications: superconductivity, electron gas in semiconductor, superfluidity of the liquid helium.' 'Macroscopic and microscopic systems. Postulates of thermodynamics. Fundamental equation. Equations of state. Thermodynamic equilibrium. Thermodynamic derivatives. Thermodynamic potentials. Maxwell relations. Born diagram. Reduction of thermodynamic derivatives.^lMicrocanonical formalism. Boltzmann equation. Einstein model of a crystalline solid. Canonical formalism. Example: ideal classical gas. Maxwell-Boltzmann distribution. Grand canonical formalism. Quantum gases. Fermions and bosons. Bose-Einstein distribution. Fermi-Dirac distribution. Examples: electron gas and photon gas. Planck distribution.^lApplications: superconductivity, electron gas in semiconductor, superfluidity of the liquid helium.' | Out-Null

# 3) Bibliografia
Replace-WithLineBreaks 'CALLEN, H.B., Thermodynamics and an introduction to thermostatistics, John Wiley & Sons, New York, 1985. SALINAS, S. R. A., Introdução à Física Estatística, Edusp, São Paulo, 1999. CASQUILHO J.P. e TEIXEIRA P.I.C, Introdução à Física Estatística, Editora Livraria da Física, São Paulo, 2012. DALARSSON, N. DALARSSON, M. GOLUBOVIC, L.  Introductory Statistical Thermodynamics. Academic Press, 2011.' 'CALLEN, H.B., Thermodynamics and an introduction to thermostatistics, John Wiley & Sons, New York, 1985. ^lSALINAS, S. R. A., Introdução à Física Estatística, Edusp, São Paulo, 1999. ^lCASQUILHO J.P. e TEIXEIRA P.I.C, Introdução à Física Estatística, Editora Livraria da Física, São Paulo, 2012. ^lDALARSSON, N. DALARSSON, M. GOLUBOVIC, L.  Introductory Statistical Thermodynamics. Academic Press, 2011.' | Out-Null
